# Fruta / hortaliza, semanal
# Two new weekly price records were added to the top of the data table
# (rows 9 and 10), pushing all the existing data rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the data block (rows 9 and 10),
# shifting all existing data rows (previously 9..84) down to 11..86.
$ws.Rows("9:10").Insert()

# Row 9: new "Especial" quality record
$ws.Cells.Item(9, 1).Value2 = 5
$ws.Cells.Item(9, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(9, 3).Value2 = "Maule"
$ws.Cells.Item(9, 4).Value2 = 45050
$ws.Cells.Item(9, 5).Value2 = 7
$ws.Cells.Item(9, 6).Value2 = "Fruta"
$ws.Cells.Item(9, 7).Value2 = 100104
$ws.Cells.Item(9, 8).Value2 = "Frutos de pepita"
$ws.Cells.Item(9, 9).Value2 = 100104003
$ws.Cells.Item(9, 10).Value2 = "Membrillo"
$ws.Cells.Item(9, 11).Value2 = "Champion"
$ws.Cells.Item(9, 12).Value2 = "Especial"
$ws.Cells.Item(9, 13).Value2 = 180
$ws.Cells.Item(9, 14).Value2 = 12000
$ws.Cells.Item(9, 15).Value2 = 12000
$ws.Cells.Item(9, 16).Value2 = 12000
$ws.Cells.Item(9, 17).Value2 = "$/caja 18 kilos granel"
$ws.Cells.Item(9, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(9, 19).Value2 = 667
$ws.Cells.Item(9, 20).Value2 = 18

# Row 10: new "Primera" quality record
$ws.Cells.Item(10, 1).Value2 = 5
$ws.Cells.Item(10, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(10, 3).Value2 = "Maule"
$ws.Cells.Item(10, 4).Value2 = 45050
$ws.Cells.Item(10, 5).Value2 = 7
$ws.Cells.Item(10, 6).Value2 = "Fruta"
$ws.Cells.Item(10, 7).Value2 = 100104
$ws.Cells.Item(10, 8).Value2 = "Frutos de pepita"
$ws.Cells.Item(10, 9).Value2 = 100104003
$ws.Cells.Item(10, 10).Value2 = "Membrillo"
$ws.Cells.Item(10, 11).Value2 = "Champion"
$ws.Cells.Item(10, 12).Value2 = "Primera"
$ws.Cells.Item(10, 13).Value2 = 210
$ws.Cells.Item(10, 14).Value2 = 10000
$ws.Cells.Item(10, 15).Value2 = 10000
$ws.Cells.Item(10, 16).Value2 = 10000
$ws.Cells.Item(10, 17).Value2 = "$/caja 18 kilos granel"
$ws.Cells.Item(10, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(10, 19).Value2 = 556
$ws.Cells.Item(10, 20).Value2 = 18
